{"js": "// Replace the 25 three-digit-by-one-digit multiplication equations with\n// their updated values (old -> new), matching the commit's canonical diff.\nconst replacements = [\n  [\"305\u00d77=2135\", \"452\u00d78=3616\"],\n  [\"114\u00d76=684\", \"346\u00d79=3114\"],\n  [\"264\u00d78=2112\", \"390\u00d78=3120\"],\n  [\"816\u00d74=3264\", \"972\u00d75=4860\"],\n  [\"312\u00d73=936\", \"512\u00d73=1536\"],\n  [\"482\u00d73=1446\", \"636\u00d74=2544\"],\n  [\"944\u00d74=3776\", \"731\u00d76=4386\"],\n  [\"586\u00d77=4102\", \"465\u00d75=2325\"],\n  [\"464\u00d74=1856\", \"351\u00d75=1755\"],\n  [\"907\u00d76=5442\", \"164\u00d78=1312\"],\n  [\"951\u00d75=4755\", \"730\u00d76=4380\"],\n  [\"428\u00d79=3852\", \"293\u00d72=586\"],\n  [\"304\u00d74=1216\", \"965\u00d79=8685\"],\n  [\"710\u00d78=5680\", \"472\u00d76=2832\"],\n  [\"180\u00d75=900\", \"389\u00d76=2334\"],\n  [\"129\u00d72=258\", \"531\u00d77=3717\"],\n  [\"828\u00d79=7452\", \"926\u00d77=6482\"],\n  [\"596\u00d73=1788\", \"785\u00d75=3925\"],\n  [\"250\u00d72=500\", \"263\u00d78=2104\"],\n  [\"859\u00d76=5154\", \"436\u00d79=3924\"],\n  [\"493\u00d73=1479\", \"742\u00d75=3710\"],\n  [\"928\u00d75=4640\", \"582\u00d75=2910\"],\n  [\"190\u00d77=1330\", \"767\u00d75=3835\"],\n  [\"795\u00d77=5565\", \"471\u00d79=4239\"],\n  [\"836\u00d77=5852\", \"907\u00d79=8163\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication equations with\n# their updated values (old -> new), matching the commit's canonical diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"305\u00d77=2135\", \"452\u00d78=3616\"),\n    @(\"114\u00d76=684\", \"346\u00d79=3114\"),\n    @(\"264\u00d78=2112\", \"390\u00d78=3120\"),\n    @(\"816\u00d74=3264\", \"972\u00d75=4860\"),\n    @(\"312\u00d73=936\", \"512\u00d73=1536\"),\n    @(\"482\u00d73=1446\", \"636\u00d74=2544\"),\n    @(\"944\u00d74=3776\", \"731\u00d76=4386\"),\n    @(\"586\u00d77=4102\", \"465\u00d75=2325\"),\n    @(\"464\u00d74=1856\", \"351\u00d75=1755\"),\n    @(\"907\u00d76=5442\", \"164\u00d78=1312\"),\n    @(\"951\u00d75=4755\", \"730\u00d76=4380\"),\n    @(\"428\u00d79=3852\", \"293\u00d72=586\"),\n    @(\"304\u00d74=1216\", \"965\u00d79=8685\"),\n    @(\"710\u00d78=5680\", \"472\u00d76=2832\"),\n    @(\"180\u00d75=900\", \"389\u00d76=2334\"),\n    @(\"129\u00d72=258\", \"531\u00d77=3717\"),\n    @(\"828\u00d79=7452\", \"926\u00d77=6482\"),\n    @(\"596\u00d73=1788\", \"785\u00d75=3925\"),\n    @(\"250\u00d72=500\", \"263\u00d78=2104\"),\n    @(\"859\u00d76=5154\", \"436\u00d79=3924\"),\n    @(\"493\u00d73=1479\", \"742\u00d75=3710\"),\n    @(\"928\u00d75=4640\", \"582\u00d75=2910\"),\n    @(\"190\u00d77=1330\", \"767\u00d75=3835\"),\n    @(\"795\u00d77=5565\", \"471\u00d79=4239\"),\n    @(\"836\u00d77=5852\", \"907\u00d79=8163\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
